$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the activity codes (CPS1.1 -> CPS1.3, CPS1.2 -> CPS1.4)
$ws.Range("B2").Value = "CPS1.3"
$ws.Range("B3").Value = "CPS1.4"

# Copy the style used on B2 onto B3
$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore B3's text value (paste special formats may have touched it)
$ws.Range("B3").Value = "CPS1.4"

# Update the active cell / selection
$ws.Range("B4").Select()
